$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("README")

# Insert three new rows above the old row 9 ("Reference:") -- these become
# rows 9, 10 and 11; everything from the old row 9 onward shifts down by 3.
$ws2.Rows("9:11").Insert()

# New row 9: a plain paragraph of explanatory text (same look as row 8 / row 4 / row 6).
$ws2.Range("A8:B8").Copy()
$ws2.Range("A9:B9").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("A9").RowHeight = 30

# New row 10: the dataset URL, turned into a real hyperlink (add the
# hyperlink before re-applying the B7-style formatting so the final cell
# keeps the same shared hyperlink style as the existing VASCAN link instead
# of the ad-hoc style Hyperlinks.Add() would otherwise stamp on it).
$ws2.Range("B10").Value2 = "http://www.gbif.org/dataset/3f8a1297-3259-4700-91fc-acc4170b27ce"
$ws2.Hyperlinks.Add($ws2.Range("B10"), "http://www.gbif.org/dataset/3f8a1297-3259-4700-91fc-acc4170b27ce")

# New rows 10-11: hyperlink-styled cells (same look as row 7, which already
# carries the existing VASCAN hyperlink formatting).
$ws2.Range("A7:B7").Copy()
$ws2.Range("A10:B11").PasteSpecial(-4122) # xlPasteFormats

# Fill in the remaining new cell value (descriptive text). Written after the
# URL so the new shared-string entries land in the same order as the source
# file (URL first, description second).
$ws2.Range("B9").Value2 = "The Database of Vascular Plants of Canada is used in building the GBIF Backbone Taxonomy. To see how GBIF indexes the data go here:"
$ws2.Range("B11").Value2 = ""

# Make README the active tab (was Classification before).
$ws2.Activate()
